$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "glossary" -- update header row and fill in newly required
# required/field_type values for each data-dictionary row.
# ---------------------------------------------------------------------------
$glossary = $wb.Worksheets.Item(1)

# Header row
$glossary.Range("C1").Value = "required"
$glossary.Range("D1").Value = "field_type"
$glossary.Range("E1").Value = "unit"
$glossary.Range("F1").Value = "variable_codes"

# study_id row
$glossary.Range("C2").Value = "optional"
$glossary.Range("D2").Value = "character"

# site_id row
$glossary.Range("C3").Value = "optional"
$glossary.Range("D3").Value = "character"

# core_id row
$glossary.Range("C4").Value = "optional"
$glossary.Range("D4").Value = "character"

# impact_class row
$glossary.Range("C5").Value = "optional"
$glossary.Range("D5").Value = "factor"

# ---------------------------------------------------------------------------
# Sheet 2: "impacts" -- trim to the four required id/class columns, drop the
# example data row, and size the columns for the new layout.
# ---------------------------------------------------------------------------
$impacts = $wb.Worksheets.Item(2)

$impacts.Range("A1").Value = "study_id"
$impacts.Range("B1").Value = "site_id"
$impacts.Range("C1").Value = "core_id"
$impacts.Range("D1").Value = "impact_class"

# Remove the old example-value row and the now-unused E column entirely.
$impacts.Range("A2:E2").ClearContents()
$impacts.Columns.Item(5).Delete()

$impacts.Columns.Item(1).ColumnWidth = 7.83
$impacts.Columns.Item(2).ColumnWidth = 6.83
$impacts.Columns.Item(3).ColumnWidth = 6.83
$impacts.Columns.Item(4).ColumnWidth = 11.83
